# StructureDefinition-aliquot-concentration.xlsx — gh-pages IG regeneration
#
# The underlying IG build ran again (new publish Date, and this copy of the
# spreadsheet was regenerated against FHIR R4 instead of R4B), which changes
# a handful of values on the "Metadata" sheet plus a few StructureDefinition
# snapshot cells on the "Elements" sheet (root ele-1 constraint text loses the
# "unless an empty Parameters resource ... or $this is Parameters" carve-out,
# Extension.id's declared type flips from "id" to "string", and the
# Extensibility doc link used in Extension.value[x]'s definition moves from
# /R4B/ to /R4/).

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ---------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Date
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# FHIR Version
$wsMeta.Range("B15").Value = "4.0.1"

# ---- Elements sheet -----------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Row 2 ("Extension"): Constraint(s) drops the Parameters-resource carve-out
# from the ele-1 invariant, matching the text already used further down the
# snapshot (e.g. Extension.extension's ele-1).
$wsElem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 ("Extension.id"): Type(s) changes from "id" to "string".
$wsElem.Range("K3").Value = "string
"

# Row 6 ("Extension.value[x]"): Definition's Extensibility link moves from
# the R4B docs to the R4 docs.
$wsElem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
